# Q3 Update - 2025
# Applies the changes described in the diff:
#  - Renames the shared short-url string "K5jvPY" -> "H3gNzy" for every row that uses it
#  - Updates Cuba/2024 row (row 28): refugees 49 -> 50, asylum_seekers 93 -> 15
#  - Updates Haiti/2024 row (row 29): ooc 78 -> 32
#  - Updates Venezuela/2024 row (row 30): ooc 144 -> 58

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row and rewrite every matching short-url cell in column B.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "K5jvPY") {
        $cell.Value = "H3gNzy"
    }
}

# The source sheet stores every numeric-looking field as text (shared
# strings), so force text formatting before writing these new values to
# keep them as text instead of letting them be parsed back into numbers.

# Row 28 = Cuba, year 2024: refugees 49 -> 50, asylum_seekers 93 -> 15
$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = "50"
$ws.Range("O28").NumberFormat = "@"
$ws.Range("O28").Value = "15"

# Row 29 = Haiti, year 2024: ooc (others of concern) 78 -> 32
$ws.Range("T29").NumberFormat = "@"
$ws.Range("T29").Value = "32"

# Row 30 = Venezuela, year 2024: ooc (others of concern) 144 -> 58
$ws.Range("T30").NumberFormat = "@"
$ws.Range("T30").Value = "58"
